$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlLeft = -4131
$dateFmt = "[$-409]dd\-mmm\-yy;@"

# ---------------------------------------------------------------------------
# Move the existing Windows 10 test-matrix rows (13-15) down to rows 17-19
# (leaving row 16 blank), before overwriting 13-15 with the new Windows 11
# rows. Values are written explicitly (rather than copy/paste) so styles are
# reproduced deterministically via the same formatting calls used below.
# ---------------------------------------------------------------------------

function Set-MatrixRow {
    param($row, $a, $b, $c, $d, $e, $f, $g, $h)

    $ws.Range("A$row").Value = $a
    $ws.Range("B$row").Value = $b
    $ws.Range("B$row").HorizontalAlignment = $xlCenter
    $ws.Range("C$row").Value = $c
    $ws.Range("C$row").HorizontalAlignment = $xlCenter
    $ws.Range("D$row").Value = $d
    $ws.Range("D$row").HorizontalAlignment = $xlCenter
    $ws.Range("E$row").Value = $e
    $ws.Range("E$row").HorizontalAlignment = $xlLeft
    $ws.Range("F$row").Value = $f
    $ws.Range("F$row").HorizontalAlignment = $xlCenter
    $ws.Range("G$row").Value = $g
    $ws.Range("G$row").HorizontalAlignment = $xlCenter
    $ws.Range("G$row").NumberFormat = $dateFmt
    $ws.Range("H$row").Value = $h
    $ws.Range("H$row").HorizontalAlignment = $xlCenter
}

# Rows 17-19: the old Windows 10 rows, relocated from 13-15.
Set-MatrixRow 17 "Windows 10" "Pro 21H2" "x86_64" "8 GB" "VS2017 15.9.50" "VMware" 44824 "PowerShell 7.2.6"
Set-MatrixRow 18 "Windows 10" "Pro 21H2" "x86_64" "8 GB" "VS2019 16.11.19" "VMware" 44824 "PowerShell 7.2.6"
Set-MatrixRow 19 "Windows 10" "Pro 21H2" "x86_64" "8 GB" "VS2022 17.3.4" "VMware" 44824 "PowerShell 7.2.6"

# Rows 13-15: new Windows 11 rows.
# New shared strings must be introduced in this exact order so the
# sharedStrings.xml table indices line up: Windows 11, Pro 22H2,
# PowerShell 7.1.3, VS2019 16.11.21, VS2017 15.9.51, VS2022 17.4.1.
$ws.Range("A13").Value = "Windows 11"
$ws.Range("B13").Value = "Pro 22H2"
$ws.Range("B13").HorizontalAlignment = $xlCenter
$ws.Range("H13").Value = "PowerShell 7.1.3"
$ws.Range("H13").HorizontalAlignment = $xlCenter
$ws.Range("E14").Value = "VS2019 16.11.21"
$ws.Range("E14").HorizontalAlignment = $xlLeft
$ws.Range("E13").Value = "VS2017 15.9.51"
$ws.Range("E13").HorizontalAlignment = $xlLeft
$ws.Range("E15").Value = "VS2022 17.4.1"
$ws.Range("E15").HorizontalAlignment = $xlLeft

$ws.Range("C13").Value = "x86_64"
$ws.Range("C13").HorizontalAlignment = $xlCenter
$ws.Range("D13").Value = "8 GB"
$ws.Range("D13").HorizontalAlignment = $xlCenter
$ws.Range("F13").Value = "VMware"
$ws.Range("F13").HorizontalAlignment = $xlCenter
$ws.Range("G13").Value = 44888
$ws.Range("G13").HorizontalAlignment = $xlCenter
$ws.Range("G13").NumberFormat = $dateFmt

$ws.Range("A14").Value = "Windows 11"
$ws.Range("B14").Value = "Pro 21H2"
$ws.Range("B14").HorizontalAlignment = $xlCenter
$ws.Range("C14").Value = "x86_64"
$ws.Range("C14").HorizontalAlignment = $xlCenter
$ws.Range("D14").Value = "8 GB"
$ws.Range("D14").HorizontalAlignment = $xlCenter
$ws.Range("F14").Value = "VMware"
$ws.Range("F14").HorizontalAlignment = $xlCenter
$ws.Range("G14").Value = 44888
$ws.Range("G14").HorizontalAlignment = $xlCenter
$ws.Range("G14").NumberFormat = $dateFmt
$ws.Range("H14").Value = "PowerShell 7.1.3"
$ws.Range("H14").HorizontalAlignment = $xlCenter

$ws.Range("A15").Value = "Windows 11"
$ws.Range("B15").Value = "Pro 21H2"
$ws.Range("B15").HorizontalAlignment = $xlCenter
$ws.Range("C15").Value = "x86_64"
$ws.Range("C15").HorizontalAlignment = $xlCenter
$ws.Range("D15").Value = "8 GB"
$ws.Range("D15").HorizontalAlignment = $xlCenter
$ws.Range("F15").Value = "VMware"
$ws.Range("F15").HorizontalAlignment = $xlCenter
$ws.Range("G15").Value = 44888
$ws.Range("G15").HorizontalAlignment = $xlCenter
$ws.Range("G15").NumberFormat = $dateFmt
$ws.Range("H15").Value = "PowerShell 7.1.3"
$ws.Range("H15").HorizontalAlignment = $xlCenter

# Refresh the "this chart updated" timestamp formula cell.
$ws.Range("C4").Formula = "=TODAY()"

# Move the active selection to A3, matching the post-edit cursor position.
$ws.Range("A3").Select() | Out-Null
